# REPORTGEN-709: update full detailed excel reports
#
# On the "Summary" sheet of the CISQ Full Detailed Report template:
#   - the B3:C3 and B4:C4 merged header cells are split back into separate
#     cells (only the F1:G1 merge remains),
#   - the "RepGen:TEXT;APPLICATION_NAME" / "RepGen:TEXT;LAST_SNAPSHOT_DATE"
#     placeholder text that used to sit in D3 / D4 (to the right of the
#     merged B:C label cell) now sits directly in C3 / C4, right next to
#     the label, leaving D3 / D4 blank,
#   - the sheet's saved cursor/selection moves from B9 to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Split the merged label cells back into individual cells.
$ws.Range("B3:C3").UnMerge()
$ws.Range("B4:C4").UnMerge()

# Move (not copy) the placeholder values that lived in D3/D4 into the
# freshly-unmerged C3/C4 cells; Cut leaves the source cell blank while
# keeping its existing cell formatting, and the destination picks up
# both the value and the formatting of the source cell.
$ws.Range("D3").Cut($ws.Range("C3"))
$ws.Range("D4").Cut($ws.Range("C4"))

# Update the sheet's stored selection.
$ws.Range("B7").Select() | Out-Null
